# entrega final - script
# Update the "livros e artigos para tese" sheet with the new species rows
# and refresh the view state (tab ratio / selected cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("livros e artigos para tese")

# Row 2 previously held "Justia pectoralis Jacq." - replace it, then append
# three new rows with the additional macrofita species.
$ws.Range("A2").Value = "Panicum pernambuncense (Spreng.) Mez ex Pilg."
$ws.Range("A3").Value = "Reimarochloa acuta (Flüggé) Hitchc."
$ws.Range("A4").Value = "Sacciolepis myuros (Lam.) Chase"
$ws.Range("A5").Value = "Urochloa plantaginea (Link) R. D. Webster"

# Restore the workbook tab ratio (bookViews/workbookView@tabRatio="925")
$excel.ActiveWindow.TabRatio = 0.925

# Leave the selection on A11, where editing left off.
[void]$ws.Range("A11").Select()
